# Trading update: 2026-02-17 08:43:06
# Appends the newest trade (Trade # 56, row 57) to both the "All Trades"
# and "MarketMaking" sheets — a new OPEN MarketMaking position.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 57

    $ws.Cells.Item($row, 1).Value = 56

    # Dates like "2026-02-17" get auto-parsed into date serials by the
    # Value setter, so force text formatting, write, then reset the style
    # back to Normal so no stray numFmt/style survives on the cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "08:42:39"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.53
    $ws.Cells.Item($row, 7).Value = ""
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 99.56427336807459
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = ""
    $ws.Cells.Item($row, 17).Value = 0
}
